# Auto-update draw results: append the 2025-09-23 Pick 3 row.
#
# The source data (date, phase code, and numeric-looking draw result) must
# land in the sheet as literal TEXT, matching every other row already in
# the table -- not as Excel's auto-detected dates/numbers. We stage each
# value as a quoted-string formula (guaranteed text result), then collapse
# it to a plain value in place via Copy/PasteSpecial(values-only) so no
# formula and no extra number-format/style survives the write.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7

$ws.Range("A$row").Formula = "=""2025-09-23"""
$ws.Range("B$row").Formula = "=""Pick 3"""
$ws.Range("C$row").Formula = "=""250923"""
$ws.Range("D$row").Formula = "=""7-4-5"""
$ws.Range("E$row").Formula = "=""2025-09-23T21:36:20.956+04:00"""

$target = $ws.Range("A$($row):E$row")
$target.Copy()
$target.PasteSpecial(-4163)
